$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '79.788.83'
$ws.Range("E2").Value = '  +4.54%  '

# Row 3
$ws.Range("D3").Value = '3.204.27'
$ws.Range("E3").Value = '  +5.40%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '''206.19'
$ws.Range("E5").Value = '  +3.02%  '

# Row 6
$ws.Range("D6").Value = '''637.35'
$ws.Range("E6").Value = '  +2.15%  '

# Row 7
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("E8").Value = '  +19.21%  '

# Row 9
$ws.Range("D9").Value = '''0.610'
$ws.Range("E9").Value = '  +11.35%  '

# Row 10
$ws.Range("D10").Value = '3.204.19'
$ws.Range("E10").Value = '  +5.39%  '

# Row 11
$ws.Range("D11").Value = '''0.624'
$ws.Range("E11").Value = '  +42.02%  '

# Row 12
$ws.Range("D12").Value = '''0.0000252'
$ws.Range("E12").Value = '  +30.37%  '

# Row 13
$ws.Range("E13").Value = '  +3.39%  '

# Row 14
$ws.Range("D14").Value = '''5.43'
$ws.Range("E14").Value = '  +3.82%  '

# Row 15
$ws.Range("D15").Value = '3.795.01'
$ws.Range("E15").Value = '  +5.45%  '

# Row 16
$ws.Range("D16").Value = '''32.49'
$ws.Range("E16").Value = '  +11.62%  '

# Row 17
$ws.Range("D17").Value = '79.583.39'
$ws.Range("E17").Value = '  +4.34%  '

# Row 18
$ws.Range("D18").Value = '3.200.07'
$ws.Range("E18").Value = '  +5.66%  '

# Row 19
$ws.Range("D19").Value = '''14.67'
$ws.Range("E19").Value = '  +8.18%  '

# Row 20
$ws.Range("D20").Value = '''9.46'
$ws.Range("E20").Value = '  +5.51%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '''443.51'
$ws.Range("E21").Value = '  +18.33%  '

# Row 22
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").Value = '''2.95'
$ws.Range("E22").Value = '  +27.40%  '

# Row 23
$ws.Range("D23").Value = '''5.27'
$ws.Range("E23").Value = '  +20.94%  '

# Row 24
$ws.Range("D24").Value = '''4.85'
$ws.Range("E24").Value = '  +11.27%  '

# Row 25
$ws.Range("D25").Value = '''77.83'
$ws.Range("E25").Value = '  +6.38%  '

# Row 26
$ws.Range("D26").Value = '''10.95'
$ws.Range("E26").Value = '  +11.64%  '

# Row 27
$ws.Range("E27").Value = '  +0.40%  '

# Row 28
$ws.Range("E28").Value = '  +9.68%  '

# Row 29
$ws.Range("D29").Value = '''9.29'
$ws.Range("E29").Value = '  +12.03%  '

# Row 30
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  +0.04%  '

# Row 31
$ws.Range("E31").Value = '  +6.58%  '

# Row 32
$ws.Range("D32").Value = '''533.70'
$ws.Range("E32").Value = '  +8.57%  '

# Row 33
$ws.Range("E33").Value = '  +3.86%  '

# Row 34
$ws.Range("D34").Value = '''0.146'
$ws.Range("E34").Value = '  +26.90%  '

# Row 35
$ws.Range("D35").Value = '''23.38'
$ws.Range("E35").Value = '  +13.16%  '

# Row 36
$ws.Range("E36").Value = '  +18.17%  '

# Row 37
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.04%  '

# Row 38
$ws.Range("E38").Value = '  +7.06%  '

# Row 39
$ws.Range("D39").Value = '''164.78'
$ws.Range("E39").Value = '  +1.28%  '

# Row 40
$ws.Range("E40").Value = '  +0.03%  '

# Row 41
$ws.Range("D41").Value = '''192.46'
$ws.Range("E41").Value = '  +1.32%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").Value = '''5.58'
$ws.Range("E43").Value = '  +8.90%  '

# Row 44
$ws.Range("D44").Value = '''1.84'
$ws.Range("E44").Value = '  +11.56%  '

# Row 45
$ws.Range("E45").Value = '  +0.90%  '

# Row 46
$ws.Range("E46").Value = '  +4.59%  '

# Row 47
$ws.Range("D47").Value = '''43.54'
$ws.Range("E47").Value = '  +3.53%  '

# Row 48
$ws.Range("E48").Value = '  +5.98%  '

# Row 49
$ws.Range("D49").Value = '''25.80'
$ws.Range("E49").Value = '  +15.96%  '

# Row 50
$ws.Range("D50").Value = '''0.640'
$ws.Range("E50").Value = '  +5.93%  '

# Row 51
$ws.Range("E51").Value = '  +7.89%  '
